# Apply updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.748.11"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.038.52"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.607"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.64"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0838"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.77%  "

$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.339.71"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.45"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.45"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.771"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.054.91"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.735.85"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("E19").Value = "  -0.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0822"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("E22").Value = "  -0.74%  "

$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("E25").Value = "  +2.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.37"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.84%  "

$ws.Range("E28").Value = "  -1.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.79"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("E30").Value = "  -0.24%  "

$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("E32").Value = "  +8.20%  "

$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("E34").Value = "  +0.79%  "

$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.39"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("E37").Value = "  +3.53%  "

$ws.Range("E38").Value = "  +5.57%  "

$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.06"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +9.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.527.79"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.31"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.77%  "

$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.85"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.28"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.24%  "

$ws.Range("E46").Value = "  -1.30%  "

$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.00"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.228.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.56%  "
